$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column D ("Price") while forcing it to stay a
# text value (these price strings look numeric, e.g. "260.49", and Excel
# would otherwise silently convert them to floating point numbers, which
# changes both the stored representation and the cell type). Briefly mark
# the cell as Text, assign the literal string, then strip the formatting
# back off so the cell is left without any explicit style, matching the
# rest of the sheet.
function Set-PriceText($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Simple price (column D) updates ---
Set-PriceText "D2"  "260.49"
Set-PriceText "D4"  "6.183"
Set-PriceText "D5"  "0.06101"
Set-PriceText "D6"  "6.740"
Set-PriceText "D7"  "3.483"

# --- Rows 10-18: list shifted down by one (new "One" entry inserted at rank 9) ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-PriceText "D10" "0.01328"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-PriceText "D11" "0.1579"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-PriceText "D12" "0.08053"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-PriceText "D13" "0.03321"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-PriceText "D14" "0.03049"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-PriceText "D15" "0.09301"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-PriceText "D16" "3.894"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-PriceText "D17" "0.001696"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-PriceText "D18" "0.04841"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- More simple price (column D) updates ---
Set-PriceText "D19" "0.006186"
Set-PriceText "D20" "0.001101"
Set-PriceText "D21" "0.003396"
Set-PriceText "D22" "0.0001501"
Set-PriceText "D23" "3.692"
Set-PriceText "D27" "0.0003018"
Set-PriceText "D40" "0.04596"
Set-PriceText "D41" "0.007146"

# --- Rows 42-43: swap of CEJI / BKEXToken ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-PriceText "D42" "0.003903"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-PriceText "D43" "0.1118"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Remaining simple price (column D) updates ---
Set-PriceText "D44" "0.01069"
Set-PriceText "D45" "0.002974"
Set-PriceText "D46" "0.00005937"
Set-PriceText "D47" "0.00000000751"
Set-PriceText "D48" "0.7506"

# --- Row 49: price + volume label update ---
Set-PriceText "D49" "0.06432"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"

Set-PriceText "D51" "0.01011"
